# =========================================================================
# Edit: add "backup" column (R) that preserves the prior detect_structure
# values for the rows where that heuristic still applies, reset the
# detect_structure (Q) column to 0 for the rows it no longer drives, flip the
# isPivot flag on the final original row, and append six new monthly candles
# (Jul-2024 .. Dec-2024).
# =========================================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New column header R1 = "backup" (reuse header style from Q1) ---
$ws.Range("Q1").Copy() | Out-Null
$ws.Range("R1").PasteSpecial(-4122) | Out-Null
$ws.Range("R1").Value = "backup"

# --- 2. Populate column R (backup) for rows 2-343, default 0 ---
for ($i = 2; $i -le 343; $i++) {
    $ws.Cells.Item($i, 18).Value = 0
}

# --- 3. Rows whose prior detect_structure value is preserved into backup ---
$ws.Cells.Item(207, 18).Value = 1
$ws.Cells.Item(217, 18).Value = 1
$ws.Cells.Item(220, 18).Value = 2
$ws.Cells.Item(223, 18).Value = 1
$ws.Cells.Item(226, 18).Value = 2
$ws.Cells.Item(236, 18).Value = 1
$ws.Cells.Item(241, 18).Value = 1
$ws.Cells.Item(332, 18).Value = 1
$ws.Cells.Item(336, 18).Value = 2

# --- 4. detect_structure (Q) reset to 0 now that backup holds the value ---
$ws.Cells.Item(16, 17).Value = 0
$ws.Cells.Item(23, 17).Value = 0
$ws.Cells.Item(27, 17).Value = 0
$ws.Cells.Item(30, 17).Value = 0
$ws.Cells.Item(33, 17).Value = 0
$ws.Cells.Item(39, 17).Value = 0
$ws.Cells.Item(45, 17).Value = 0
$ws.Cells.Item(56, 17).Value = 0
$ws.Cells.Item(57, 17).Value = 0

# --- 5. isPivot (O) flips to 1 on the last pre-existing row (343) ---
$ws.Cells.Item(343, 15).Value = 1

# --- 6. Append six new monthly rows (344-349); column F (Adj Close) and
#        column R (backup) are left blank, matching the source rows. ---
$newA = New-Object 'object[,]' 6,5
$newG = New-Object 'object[,]' 6,11
$newA[0,0]=45474; $newA[0,1]=1557.239718276524; $newA[0,2]=1603.358215488316; $newA[0,3]=1458.05134141316; $newA[0,4]=1500.332885742188
$newA[1,0]=45505; $newA[1,1]=1506.885643294071; $newA[1,2]=1534.516838843759; $newA[1,3]=1428.402016039172; $newA[1,4]=1504.518676757812
$newA[2,0]=45536; $newA[2,1]=1510.875; $newA[2,2]=1533.474975585938; $newA[2,3]=1445.875; $newA[2,4]=1476.574951171875
$newA[3,0]=45566; $newA[3,1]=1480.650024414062; $newA[3,2]=1487.949951171875; $newA[3,3]=1320.300048828125; $newA[3,4]=1332.050048828125
$newA[4,0]=45597; $newA[4,1]=1333.050048828125; $newA[4,2]=1341.949951171875; $newA[4,3]=1217.25; $newA[4,4]=1292.199951171875
$newA[5,0]=45627; $newA[5,1]=1288; $newA[5,2]=1329.949951171875; $newA[5,3]=1201.5; $newA[5,4]=1210.699951171875
$newG[0,0]=230090166; $newG[0,1]=2024; $newG[0,2]=7; $newG[0,3]=1; $newG[0,4]=0; $newG[0,5]=0; $newG[0,6]=0; $newG[0,7]=27; $newG[0,8]=0; $newG[0,9]=0; $newG[0,10]=0
$newG[1,0]=259569538; $newG[1,1]=2024; $newG[1,2]=8; $newG[1,3]=1; $newG[1,4]=0; $newG[1,5]=0; $newG[1,6]=0; $newG[1,7]=31; $newG[1,8]=0; $newG[1,9]=0; $newG[1,10]=0
$newG[2,0]=297714904; $newG[2,1]=2024; $newG[2,2]=9; $newG[2,3]=1; $newG[2,4]=0; $newG[2,5]=0; $newG[2,6]=0; $newG[2,7]=35; $newG[2,8]=0; $newG[2,9]=0; $newG[2,10]=0
$newG[3,0]=400774438; $newG[3,1]=2024; $newG[3,2]=10; $newG[3,3]=1; $newG[3,4]=0; $newG[3,5]=0; $newG[3,6]=0; $newG[3,7]=40; $newG[3,8]=0; $newG[3,9]=0; $newG[3,10]=2
$newG[4,0]=279900722; $newG[4,1]=2024; $newG[4,2]=11; $newG[4,3]=1; $newG[4,4]=0; $newG[4,5]=0; $newG[4,6]=0; $newG[4,7]=44; $newG[4,8]=0; $newG[4,9]=0; $newG[4,10]=0
$newG[5,0]=282153932; $newG[5,1]=2024; $newG[5,2]=12; $newG[5,3]=1; $newG[5,4]=0; $newG[5,5]=0; $newG[5,6]=0; $newG[5,7]=48; $newG[5,8]=0; $newG[5,9]=0; $newG[5,10]=0
$ws.Range("A344:E349").Value = $newA
$ws.Range("G344:Q349").Value = $newG

# --- 7. Match the date/time display format used by the rest of column A ---
$ws.Range("A344:A349").NumberFormat = $ws.Range("A343").NumberFormat
